$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('C2').Value = 46073
$ws.Range('C3').Value = 46073
$ws.Range('C4').Value = 46073
$ws.Range('C5').Value = 46073
$ws.Range('C6').Value = 46073
$ws.Range('C7').Value = 46073
$ws.Range('A8').Value = 'A 2644-2026'
$ws.Range('B8').Value = 46037.62291666667
$ws.Range('C8').Value = 46073
$ws.Range('G8').Value = 0.7
$ws.Range('K8').Value = 1
$ws.Range('L8').Value = 0
$ws.Range('R8').Value = 'Flikbålmossa'
$ws.Range('S8').Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/artfynd/A 2644-2026 artfynd.xlsx", "A 2644-2026")'
$ws.Range('T8').Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/kartor/A 2644-2026 karta.png", "A 2644-2026")'
$ws.Range('V8').Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/klagomål/A 2644-2026 FSC-klagomål.docx", "A 2644-2026")'
$ws.Range('W8').Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/klagomålsmail/A 2644-2026 FSC-klagomål mail.docx", "A 2644-2026")'
$ws.Range('X8').Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/tillsyn/A 2644-2026 tillsynsbegäran.docx", "A 2644-2026")'
$ws.Range('Y8').Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/tillsynsmail/A 2644-2026 tillsynsbegäran mail.docx", "A 2644-2026")'
$ws.Range('A9').Value = 'A 7082-2024'
$ws.Range('B9').Value = 45343
$ws.Range('C9').Value = 46073
$ws.Range('F9').Value = 'Övriga Aktiebolag'
$ws.Range('G9').Value = 22.7
$ws.Range('R9').Value = 'Småvänderot'
$ws.Range('S9').Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/artfynd/A 7082-2024 artfynd.xlsx", "A 7082-2024")'
$ws.Range('T9').Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/kartor/A 7082-2024 karta.png", "A 7082-2024")'
$ws.Range('V9').Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/klagomål/A 7082-2024 FSC-klagomål.docx", "A 7082-2024")'
$ws.Range('W9').Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/klagomålsmail/A 7082-2024 FSC-klagomål mail.docx", "A 7082-2024")'
$ws.Range('X9').Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/tillsyn/A 7082-2024 tillsynsbegäran.docx", "A 7082-2024")'
$ws.Range('Y9').Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/tillsynsmail/A 7082-2024 tillsynsbegäran mail.docx", "A 7082-2024")'
$ws.Range('A10').Value = 'A 20090-2024'
$ws.Range('B10').Value = 45434
$ws.Range('C10').Value = 46073
$ws.Range('G10').Value = 2.4
$ws.Range('H10').Value = 1
$ws.Range('K10').Value = 0
$ws.Range('O10').Value = 0
$ws.Range('P10').Value = 0
$ws.Range('R10').Value = 'Hasselmus'
$ws.Range('S10').Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/artfynd/A 20090-2024 artfynd.xlsx", "A 20090-2024")'
$ws.Range('T10').Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/kartor/A 20090-2024 karta.png", "A 20090-2024")'
$ws.Range('V10').Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/klagomål/A 20090-2024 FSC-klagomål.docx", "A 20090-2024")'
$ws.Range('W10').Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/klagomålsmail/A 20090-2024 FSC-klagomål mail.docx", "A 20090-2024")'
$ws.Range('X10').Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/tillsyn/A 20090-2024 tillsynsbegäran.docx", "A 20090-2024")'
$ws.Range('Y10').Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/tillsynsmail/A 20090-2024 tillsynsbegäran mail.docx", "A 20090-2024")'
$ws.Range('A11').Value = 'A 45802-2022'
$ws.Range('B11').Value = 44844
$ws.Range('C11').Value = 46073
$ws.Range('G11').Value = 1.2
$ws.Range('H11').Value = 0
$ws.Range('L11').Value = 1
$ws.Range('O11').Value = 1
$ws.Range('P11').Value = 1
$ws.Range('R11').Value = 'Hartsticka'
$ws.Range('S11').Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/artfynd/A 45802-2022 artfynd.xlsx", "A 45802-2022")'
$ws.Range('T11').Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/kartor/A 45802-2022 karta.png", "A 45802-2022")'
$ws.Range('V11').Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/klagomål/A 45802-2022 FSC-klagomål.docx", "A 45802-2022")'
$ws.Range('W11').Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/klagomålsmail/A 45802-2022 FSC-klagomål mail.docx", "A 45802-2022")'
$ws.Range('X11').Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/tillsyn/A 45802-2022 tillsynsbegäran.docx", "A 45802-2022")'
$ws.Range('Y11').Formula = '=HYPERLINK("https://klasma.github.io/Logging_1291/tillsynsmail/A 45802-2022 tillsynsbegäran mail.docx", "A 45802-2022")'
$ws.Range('C12').Value = 46073
$ws.Range('C13').Value = 46073
$ws.Range('C14').Value = 46073
$ws.Range('C15').Value = 46073
$ws.Range('C16').Value = 46073
$ws.Range('C17').Value = 46073
$ws.Range('A18').Value = 'A 5028-2024'
$ws.Range('B18').Value = 45329
$ws.Range('C18').Value = 46073
$ws.Range('F18').Value = 'Övriga Aktiebolag'
$ws.Range('G18').Value = 4.5
$ws.Range('A19').Value = 'A 17110-2025'
$ws.Range('B19').Value = 45755
$ws.Range('C19').Value = 46073
$ws.Range('G19').Value = 0.3
$ws.Range('A20').Value = 'A 47874-2023'
$ws.Range('B20').Value = 45204
$ws.Range('C20').Value = 46073
$ws.Range('F20').Value = 'Övriga Aktiebolag'
$ws.Range('G20').Value = 1.4
$ws.Range('A21').Value = 'A 23052-2025'
$ws.Range('B21').Value = 45790.71023148148
$ws.Range('C21').Value = 46073
$ws.Range('G21').Value = 1.2
$ws.Range('A22').Value = 'A 12000-2025'
$ws.Range('B22').Value = 45728.61288194444
$ws.Range('C22').Value = 46073
$ws.Range('F22').Value = 'Övriga Aktiebolag'
$ws.Range('G22').Value = 2.4
$ws.Range('A23').Value = 'A 4193-2024'
$ws.Range('B23').Value = 45324
$ws.Range('C23').Value = 46073
$ws.Range('F23').Value = 'Övriga Aktiebolag'
$ws.Range('G23').Value = 1.9
$ws.Range('A24').Value = 'A 41157-2025'
$ws.Range('B24').Value = 45898.52972222222
$ws.Range('C24').Value = 46073
$ws.Range('G24').Value = 1.5
$ws.Range('A25').Value = 'A 49303-2025'
$ws.Range('B25').Value = 45938.55403935185
$ws.Range('C25').Value = 46073
$ws.Range('F25').Value = 'Övriga Aktiebolag'
$ws.Range('G25').Value = 4.5
$ws.Range('A26').Value = 'A 28269-2025'
$ws.Range('B26').Value = 45818
$ws.Range('C26').Value = 46073
$ws.Range('G26').Value = 6.7
$ws.Range('A27').Value = 'A 5224-2024'
$ws.Range('B27').Value = 45330
$ws.Range('C27').Value = 46073
$ws.Range('G27').Value = 18.4
$ws.Range('A28').Value = 'A 40584-2024'
$ws.Range('B28').Value = 45555.74299768519
$ws.Range('C28').Value = 46073
$ws.Range('G28').Value = 0.8
$ws.Range('A29').Value = 'A 8721-2023'
$ws.Range('B29').Value = 44978
$ws.Range('C29').Value = 46073
$ws.Range('G29').Value = 1.7
$ws.Range('A30').Value = 'A 30911-2024'
$ws.Range('B30').Value = 45498.5925
$ws.Range('C30').Value = 46073
$ws.Range('F30').Value = 'Övriga Aktiebolag'
$ws.Range('G30').Value = 2.6
$ws.Range('A31').Value = 'A 55504-2025'
$ws.Range('B31').Value = 45971.59123842593
$ws.Range('C31').Value = 46073
$ws.Range('G31').Value = 0.6
$ws.Range('A32').Value = 'A 25475-2023'
$ws.Range('B32').Value = 45089
$ws.Range('C32').Value = 46073
$ws.Range('G32').Value = 1.5
$ws.Range('A33').Value = 'A 30775-2025'
$ws.Range('B33').Value = 45831.6150462963
$ws.Range('C33').Value = 46073
$ws.Range('G33').Value = 0.8
$ws.Range('A34').Value = 'A 7072-2024'
$ws.Range('B34').Value = 45343
$ws.Range('C34').Value = 46073
$ws.Range('G34').Value = 5.1
$ws.Range('A35').Value = 'A 32488-2025'
$ws.Range('B35').Value = 45838
$ws.Range('C35').Value = 46073
$ws.Range('G35').Value = 4.7
$ws.Range('A36').Value = 'A 14994-2022'
$ws.Range('B36').Value = 44657.53449074074
$ws.Range('C36').Value = 46073
$ws.Range('F36').Value = 'Övriga Aktiebolag'
$ws.Range('G36').Value = 4.4
$ws.Range('A37').Value = 'A 18856-2024'
$ws.Range('B37').Value = 45426
$ws.Range('C37').Value = 46073
$ws.Range('F37').Value = 'Övriga Aktiebolag'
$ws.Range('G37').Value = 3.3
$ws.Range('A38').Value = 'A 60059-2025'
$ws.Range('B38').Value = 45993.65428240741
$ws.Range('C38').Value = 46073
$ws.Range('F38').Value = 'Övriga Aktiebolag'
$ws.Range('G38').Value = 7.8
$ws.Range('A39').Value = 'A 6830-2025'
$ws.Range('B39').Value = 45700
$ws.Range('C39').Value = 46073
$ws.Range('G39').Value = 1.7
$ws.Range('A40').Value = 'A 3392-2024'
$ws.Range('B40').Value = 45318
$ws.Range('C40').Value = 46073
$ws.Range('G40').Value = 2.7
$ws.Range('A41').Value = 'A 62413-2022'
$ws.Range('B41').Value = 44923
$ws.Range('C41').Value = 46073
$ws.Range('G41').Value = 2
$ws.Range('A42').Value = 'A 37189-2025'
$ws.Range('B42').Value = 45875
$ws.Range('C42').Value = 46073
$ws.Range('G42').Value = 1.2
$ws.Range('A43').Value = 'A 37052-2025'
$ws.Range('B43').Value = 45875
$ws.Range('C43').Value = 46073
$ws.Range('G43').Value = 0.9
$ws.Range('C44').Value = 46073
$ws.Range('A45').Value = 'A 62316-2025'
$ws.Range('B45').Value = 46006.69094907407
$ws.Range('C45').Value = 46073
$ws.Range('G45').Value = 3.5
$ws.Range('A46').Value = 'A 14000-2025'
$ws.Range('B46').Value = 45740.26850694444
$ws.Range('C46').Value = 46073
$ws.Range('G46').Value = 0.5
$ws.Range('A47').Value = 'A 62357-2025'
$ws.Range('B47').Value = 46007
$ws.Range('C47').Value = 46073
$ws.Range('F47').Value = 'Övriga Aktiebolag'
$ws.Range('G47').Value = 5.5
$ws.Range('A48').Value = 'A 18682-2023'
$ws.Range('B48').Value = 45043.60021990741
$ws.Range('C48').Value = 46073
$ws.Range('G48').Value = 2.6
$ws.Range('A49').Value = 'A 3367-2024'
$ws.Range('B49').Value = 45317
$ws.Range('C49').Value = 46073
$ws.Range('G49').Value = 0.5
$ws.Range('A50').Value = 'A 38846-2024'
$ws.Range('B50').Value = 45547.60444444444
$ws.Range('C50').Value = 46073
$ws.Range('G50').Value = 1.4
$ws.Range('A51').Value = 'A 19435-2024'
$ws.Range('B51').Value = 45429
$ws.Range('C51').Value = 46073
$ws.Range('G51').Value = 3.7
$ws.Range('A52').Value = 'A 14104-2024'
$ws.Range('B52').Value = 45392.61707175926
$ws.Range('C52').Value = 46073
$ws.Range('G52').Value = 5.3
$ws.Range('A53').Value = 'A 14106-2024'
$ws.Range('B53').Value = 45392
$ws.Range('C53').Value = 46073
$ws.Range('G53').Value = 0.6
$ws.Range('A54').Value = 'A 6006-2024'
$ws.Range('B54').Value = 45336
$ws.Range('C54').Value = 46073
$ws.Range('G54').Value = 2.9
$ws.Range('A55').Value = 'A 3393-2024'
$ws.Range('B55').Value = 45318
$ws.Range('C55').Value = 46073
$ws.Range('A56').Value = 'A 3395-2024'
$ws.Range('B56').Value = 45318
$ws.Range('C56').Value = 46073
$ws.Range('G56').Value = 2.4
$ws.Range('A57').Value = 'A 3398-2024'
$ws.Range('B57').Value = 45318
$ws.Range('C57').Value = 46073
$ws.Range('G57').Value = 3.4
$ws.Range('A58').Value = 'A 1910-2026'
$ws.Range('B58').Value = 46035
$ws.Range('C58').Value = 46073
$ws.Range('G58').Value = 0.7
$ws.Range('A59').Value = 'A 7584-2024'
$ws.Range('B59').Value = 45348
$ws.Range('C59').Value = 46073
$ws.Range('G59').Value = 6.9
$ws.Range('A60').Value = 'A 35692-2024'
$ws.Range('B60').Value = 45532
$ws.Range('C60').Value = 46073
$ws.Range('G60').Value = 0
$ws.Range('A61').Value = 'A 3402-2024'
$ws.Range('B61').Value = 45318
$ws.Range('C61').Value = 46073
$ws.Range('G61').Value = 0.5
$ws.Range('A62').Value = 'A 19439-2024'
$ws.Range('B62').Value = 45429
$ws.Range('C62').Value = 46073
$ws.Range('G62').Value = 3.1
$ws.Range('A63').Value = 'A 39238-2023'
$ws.Range('B63').Value = 45162
$ws.Range('C63').Value = 46073
$ws.Range('G63').Value = 4.2
$ws.Range('A64').Value = 'A 42374-2023'
$ws.Range('B64').Value = 45180.5609375
$ws.Range('C64').Value = 46073
$ws.Range('G64').Value = 0.5
$ws.Range('A65').Value = 'A 50134-2024'
$ws.Range('B65').Value = 45600.44069444444
$ws.Range('C65').Value = 46073
$ws.Range('G65').Value = 1.8
$ws.Range('A66').Value = 'A 20950-2023'
$ws.Range('B66').Value = 45061
$ws.Range('C66').Value = 46073
$ws.Range('G66').Value = 6
$ws.Range('A67').Value = 'A 43724-2024'
$ws.Range('B67').Value = 45569.63178240741
$ws.Range('C67').Value = 46073
$ws.Range('G67').Value = 1.7
$ws.Range('A68').Value = 'A 30667-2023'
$ws.Range('B68').Value = 45112
$ws.Range('C68').Value = 46073
$ws.Range('F68').Value = 'Övriga Aktiebolag'
$ws.Range('G68').Value = 1
$ws.Range('A69').Value = 'A 5488-2024'
$ws.Range('B69').Value = 45334
$ws.Range('C69').Value = 46073
$ws.Range('G69').Value = 0.5
$ws.Range('A70').Value = 'A 39015-2023'
$ws.Range('B70').Value = 45162
$ws.Range('C70').Value = 46073
$ws.Range('G70').Value = 9
$ws.Range('A71').Value = 'A 26471-2023'
$ws.Range('B71').Value = 45092
$ws.Range('C71').Value = 46073
$ws.Range('G71').Value = 0.4
$ws.Range('A72').Value = 'A 56764-2024'
$ws.Range('B72').Value = 45628
$ws.Range('C72').Value = 46073
$ws.Range('G72').Value = 7.5
$ws.Range('A73').Value = 'A 17102-2025'
$ws.Range('B73').Value = 45755
$ws.Range('C73').Value = 46073
$ws.Range('A74').Value = 'A 63464-2023'
$ws.Range('B74').Value = 45273
$ws.Range('C74').Value = 46073
$ws.Range('G74').Value = 0.7
$ws.Range('A75').Value = 'A 23287-2022'
$ws.Range('B75').Value = 44720
$ws.Range('C75').Value = 46073
$ws.Range('G75').Value = 3.5
$ws.Range('A76').Value = 'A 3404-2024'
$ws.Range('B76').Value = 45318
$ws.Range('C76').Value = 46073
$ws.Range('G76').Value = 1.9
$ws.Range('A77').Value = 'A 15039-2024'
$ws.Range('B77').Value = 45399
$ws.Range('C77').Value = 46073
$ws.Range('G77').Value = 1.3
$ws.Range('A78').Value = 'A 14702-2024'
$ws.Range('B78').Value = 45397.50576388889
$ws.Range('C78').Value = 46073
$ws.Range('G78').Value = 1.1
$ws.Range('A79').Value = 'A 7716-2024'
$ws.Range('B79').Value = 45349
$ws.Range('C79').Value = 46073
$ws.Range('G79').Value = 1.8
$ws.Range('A80').Value = 'A 39221-2023'
$ws.Range('B80').Value = 45162
$ws.Range('C80').Value = 46073
$ws.Range('G80').Value = 3.1
$ws.Range('A81').Value = 'A 7083-2024'
$ws.Range('B81').Value = 45343
$ws.Range('C81').Value = 46073
$ws.Range('F81').Value = 'Övriga Aktiebolag'
$ws.Range('G81').Value = 10.6
$ws.Range('A82').Value = 'A 3602-2024'
$ws.Range('B82').Value = 45320
$ws.Range('C82').Value = 46073
$ws.Range('G82').Value = 0.7
$ws.Range('A83').Value = 'A 61876-2024'
$ws.Range('B83').Value = 45652.38005787037
$ws.Range('C83').Value = 46073
$ws.Range('G83').Value = 1.4
$ws.Range('A84').Value = 'A 3396-2024'
$ws.Range('C84').Value = 46073
$ws.Range('G84').Value = 4.4
$ws.Range('A85').Value = 'A 3405-2024'
$ws.Range('B85').Value = 45318
$ws.Range('C85').Value = 46073
$ws.Range('G85').Value = 0.9
$ws.Range('A86').Value = 'A 39834-2021'
$ws.Range('B86').Value = 44417
$ws.Range('C86').Value = 46073
$ws.Range('G86').Value = 1.3
$ws.Range('A87').Value = 'A 56820-2024'
$ws.Range('B87').Value = 45628
$ws.Range('C87').Value = 46073
$ws.Range('G87').Value = 1.5
$ws.Range('A88').Value = 'A 58985-2022'
$ws.Range('B88').Value = 44903.69670138889
$ws.Range('C88').Value = 46073
$ws.Range('G88').Value = 2.6
$ws.Range('A89').Value = 'A 62149-2022'
$ws.Range('B89').Value = 44922
$ws.Range('C89').Value = 46073
$ws.Range('F89').Value = 'Övriga Aktiebolag'
$ws.Range('G89').Value = 2.9
$ws.Range('A90').Value = 'A 19693-2025'
$ws.Range('B90').Value = 45771.25997685185
$ws.Range('C90').Value = 46073
$ws.Range('G90').Value = 1.7
$ws.Range('A91').Value = 'A 53978-2022'
$ws.Range('B91').Value = 44876
$ws.Range('C91').Value = 46073
$ws.Range('G91').Value = 1.7
$ws.Range('C92').Value = 46073
$ws.Range('A93').Value = 'A 3604-2024'
$ws.Range('B93').Value = 45320
$ws.Range('C93').Value = 46073
$ws.Range('G93').Value = 4.1
$ws.Range('A94').Value = 'A 2554-2023'
$ws.Range('B94').Value = 44943.67172453704
$ws.Range('C94').Value = 46073
$ws.Range('G94').Value = 0.1
$ws.Range('A95').Value = 'A 43851-2024'
$ws.Range('B95').Value = 45572.35173611111
$ws.Range('C95').Value = 46073
$ws.Range('G95').Value = 0.7
$ws.Range('A96').Value = 'A 42016-2022'
$ws.Range('B96').Value = 44827
$ws.Range('C96').Value = 46073
$ws.Range('G96').Value = 0.8
$ws.Range('A97').Value = 'A 27426-2024'
$ws.Range('B97').Value = 45474.39143518519
$ws.Range('C97').Value = 46073
$ws.Range('G97').Value = 0.6
$ws.Range('A98').Value = 'A 45946-2022'
$ws.Range('B98').Value = 44845
$ws.Range('C98').Value = 46073
$ws.Range('G98').Value = 14
$ws.Range('A99').Value = 'A 42231-2023'
$ws.Range('B99').Value = 45180.31387731482
$ws.Range('C99').Value = 46073
$ws.Range('G99').Value = 2.4
$ws.Range('A100').Value = 'A 19137-2025'
$ws.Range('B100').Value = 45769
$ws.Range('C100').Value = 46073
$ws.Range('G100').Value = 2.2
$ws.Range('A101').Value = 'A 19141-2025'
$ws.Range('B101').Value = 45769
$ws.Range('C101').Value = 46073
$ws.Range('G101').Value = 0.8
$ws.Range('A102').Value = 'A 14103-2024'
$ws.Range('B102').Value = 45392
$ws.Range('C102').Value = 46073
$ws.Range('G102').Value = 0.4
$ws.Range('A103').Value = 'A 22194-2023'
$ws.Range('B103').Value = 45069
$ws.Range('C103').Value = 46073
$ws.Range('F103').Value = 'Övriga Aktiebolag'
$ws.Range('G103').Value = 11.4
$ws.Range('A104').Value = 'A 16258-2024'
$ws.Range('B104').Value = 45407
$ws.Range('C104').Value = 46073
$ws.Range('G104').Value = 0.5
$ws.Range('A105').Value = 'A 50138-2024'
$ws.Range('B105').Value = 45600.44368055555
$ws.Range('C105').Value = 46073
$ws.Range('G105').Value = 1.4
$ws.Range('F10').ClearContents()
$ws.Range('F19').ClearContents()
$ws.Range('F24').ClearContents()
$ws.Range('F26').ClearContents()
$ws.Range('F43').ClearContents()
$ws.Range('F46').ClearContents()
$ws.Range('F58').ClearContents()
$ws.Range('F60').ClearContents()
$ws.Range('F64').ClearContents()
$ws.Range('F76').ClearContents()
$ws.Range('F80').ClearContents()
$ws.Range('F83').ClearContents()
$ws.Range('F94').ClearContents()
$ws.Range('F96').ClearContents()
$ws.Range('F98').ClearContents()
